$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" (summary) sheet: shift the old rows down and add the
#        new 2022-Q4 entry at the top of the data, 2021-Q1 becomes a new row 4.
$wsTotal = $wb.Worksheets.Item(1)

# Clone A3's formatting onto the new A4 cell (keeps the same cell style index
# as the other data rows) before filling in the values.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q1"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0

$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.01

$wsTotal.Range("B2").Value = "2022-Q4"

# --- 2. Insert the new "2022-Q4" sheet. Duplicate the existing "2022-Q3"
#        sheet (same column layout/styling) and place the copy right before
#        it, then rename the copy and overwrite it with the Q4 figures.
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)

$wsQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$wsQ4.Name = "2022-Q4"

$wsQ4.Range("D2").Value = "2.55"
$wsQ4.Range("E2").Value = "97.28"
$wsQ4.Range("F2").Value = "0.47"
$wsQ4.Range("G2").Value = "0.0120"
$wsQ4.Range("H2").Value = 5
